$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" column (C) for rows 2 through 13 from serial date
# 45174 (2023-09-05) to 45175 (2023-09-06), keeping everything else intact.
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
